# Applies the "add souls structure, animation and etc" update:
#  - TODO Before 0.0.1: mark rows 7 and 11 ("create enemy mp, hp, xp system AND
#    STAMINA" / "enrich enemies ...") as done, with a completion date.
#  - Logs: append two new dev-log entries (2024-06-22 and 2024-06-24).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "TODO Before 0.0.1"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TODO Before 0.0.1")

# Row 7: status -> done, Done at -> 2024-06-23 (serial 45466)
$ws1.Range("C7").Value = "done"
$ws1.Range("D2").Copy()
$ws1.Range("D7").PasteSpecial(-4122)
$ws1.Range("D7").Value = 45466

# Row 11: status -> done, Done at -> 2024-06-23 (serial 45466)
$ws1.Range("C11").Value = "done"
$ws1.Range("D2").Copy()
$ws1.Range("D11").PasteSpecial(-4122)
$ws1.Range("D11").Value = 45466

$excel.CutCopyMode = $false

$ws1.Activate()
$ws1.Range("C28").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Logs"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Logs")

# Row 38: 2024-06-22
$ws3.Range("A37").Copy()
$ws3.Range("A38").PasteSpecial(-4122)
$ws3.Range("A38").Value = 45465
$ws3.Range("B38").Value = "work on glow effect (looks nice), fix lags when full screeen - resolution wsa to high - found place where it can be changed, healing animation"

# Row 39: 2024-06-24
$ws3.Range("A37").Copy()
$ws3.Range("A39").PasteSpecial(-4122)
$ws3.Range("A39").Value = 45467
$ws3.Range("B39").Value = "add souls, souls drop and collecting. Lot of fun, but source consumable… Now enemy on death drop souls, and player can collect them"

$excel.CutCopyMode = $false

$ws3.Activate()
$ws3.Range("B40").Select() | Out-Null
